$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update B93 with the new trim name (new shared string will be created)
$ws.Range("B93").Value = "RX 350 AWD F SPORT BLACK LINE SPECIAL EDITION"

# Widen column B to fit the longer text (stored xlsx width of 55 chars)
$ws.Columns.Item(2).ColumnWidth = 54.166666666666664

# Update the view: scroll position and active selection/cell
$ws.Activate()
$ws.Range("B93").Select()
$excel.ActiveWindow.ScrollRow = 65
$excel.ActiveWindow.ScrollColumn = 1
